# =========================================================================
# Kaman new UI - header & footer changes
# Sheet1 (TC30_Search_Typeahead): insert CLICK_JS / extra WAIT steps
# Sheet2 (Testdata): append EleType1 / EleType2 rows
# =========================================================================

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# -------------------------------------------------------------------------
# Sheet1: expand rows 3..6 into rows 3..11
# -------------------------------------------------------------------------

# 1) New row 3: CLICK_JS / EnableCertificate_MoreInfo / JS_ID / EleType1
$ws1.Range("A3:E3").Insert()
$ws1.Range("B3").Value = "CLICK_JS"
$ws1.Range("C3").Value = "EnableCertificate_MoreInfo"
$ws1.Range("D3").Value = "JS_ID"
$ws1.Range("E3").Value = "EleType1"
$ws1.Range("A3:E3").Borders.LineStyle = 1
$ws1.Rows.Item(3).RowHeight = 15.75

# 2) Old "WAIT" row is now row 4 -- adjust its style/height to match new layout
$ws1.Range("B4").Borders.LineStyle = 1
$ws1.Rows.Item(4).RowHeight = 15.75

# 3) New row 5: CLICK_JS / EnableCertificate_GoTOPage / JS_ID / EleType2
$ws1.Range("A5:E5").Insert()
$ws1.Range("B5").Value = "CLICK_JS"
$ws1.Range("C5").Value = "EnableCertificate_GoTOPage"
$ws1.Range("D5").Value = "JS_ID"
$ws1.Range("E5").Value = "EleType2"
$ws1.Range("A5:E5").Borders.LineStyle = 1
$ws1.Rows.Item(5).RowHeight = 15.75

# 4) New rows 6 & 7: extra WAIT rows
$ws1.Range("A6:E7").Insert()
$ws1.Range("B6").Value = "WAIT"
$ws1.Range("A6:E6").Borders.LineStyle = 1
$ws1.Rows.Item(6).RowHeight = 15.75
$ws1.Range("B7").Value = "WAIT"
$ws1.Range("A7:E7").Borders.LineStyle = 1
$ws1.Rows.Item(7).RowHeight = 15.75

# Rows 8 (CLICK) & 9 (ENTERTEXT) are the old rows 4 & 5, unchanged, just shifted.

# 5) New row 10: extra WAIT row
$ws1.Range("A10:E10").Insert()
$ws1.Range("B10").Value = "WAIT"
$ws1.Range("A10:E10").Borders.LineStyle = 1

# Row 11 (VERIFY_WEBELEMENT_PRESENT) is the old row 6, unchanged, just shifted.

# -------------------------------------------------------------------------
# Sheet2 (Testdata): append EleType1 / EleType2 rows
# -------------------------------------------------------------------------
$ws2.Range("A5").Value = "EleType1"
$ws2.Range("B5").Value = "JSElement"
$ws2.Range("A6").Value = "EleType2"
$ws2.Range("B6").Value = "JSElement"
$ws2.Range("A5:B6").Borders.LineStyle = 1

[void]$ws2.Range("A5:B6").Select()

# Re-activate sheet1 last so it keeps the tabSelected flag, matching the target.
[void]$ws1.Range("B10").Select()

Write-Host "edit complete"
